$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting old B (dbExcel) -> C, old C (WebExcel) -> D
$ws.Columns.Item(2).Insert()

# New column B header + long query text (StatQuery)
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Lymphoma :: Stage 2']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Match the wrap-text style used on A2 for the new B2 cell
$ws.Range("B2").WrapText = $true

# Give the new column B the same width as column A (columns A, C & D keep their original widths untouched)
$ws.Columns.Item(2).ColumnWidth = 75

# Update selection to A2 only
$ws.Range("A2").Select()
